$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("https://azizibank.af", "NO"),
    @("https://leighton-co.com", "NO"),
    @("https://ahg.af", "NO")
)

$r = 6
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}
